# Generate Report for Handback
#
# Mirrors the "handback" step of the localization pipeline:
#  - Overview/zh-cn/de-de "Status" cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - Each language sheet's "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns get populated now that the
#    handback xliff + markdown files exist
#  - The newly-populated "Latest Target File" cells become hyperlinks to the
#    source markdown file (same target as column A's link)
#  - A couple of columns get wider to fit the new (longer) text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$md1Name = "0c4fe174-d4d5-4219-a86a-488d26943302.md"
$md1Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/410d12103bbdc5e183e7eecbfb65dbe040c24337/e2e/0c4fe174-d4d5-4219-a86a-488d26943302.md"
$md2Name = "3791bee1-aff4-47ea-9f14-6c9784f25203.md"
$md2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/410d12103bbdc5e183e7eecbfb65dbe040c24337/e2e/3791bee1-aff4-47ea-9f14-6c9784f25203.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn (E) / de-de (F) status columns for both rows
# ---------------------------------------------------------------------
$overview.Cells.Item(2, 5).Value = $statusText
$overview.Cells.Item(3, 5).Value = $statusText
$overview.Cells.Item(2, 6).Value = $statusText
$overview.Cells.Item(3, 6).Value = $statusText

# Column widths: Status columns widen to fit the longer text (best-effort;
# Excel quantizes column widths to whole pixels so we land on the closest
# achievable width to the recorded 29.9777047293527 target).
$overview.Columns.Item(5).ColumnWidth = 29.95
$overview.Columns.Item(6).ColumnWidth = 29.95

# ---------------------------------------------------------------------
# zh-cn sheet, row 2 (0c4fe174...)
# ---------------------------------------------------------------------
$zhcn.Cells.Item(2, 3).Value = $statusText
$zhcn.Cells.Item(2, 9).Value = $md1Name
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(2, 9), $md1Url, "", "", $md1Name) | Out-Null
$zhcn.Cells.Item(2, 10).Value = "0c4fe174-d4d5-4219-a86a-488d26943302.c55e58c7ab2c895dc8aced18ebd64001454b9a1e.zh-cn.xlf"
$zhcn.Cells.Item(2, 11).Value = "2016-08-17 12:59:08"

# zh-cn sheet, row 3 (3791bee1...)
$zhcn.Cells.Item(3, 3).Value = $statusText
$zhcn.Cells.Item(3, 9).Value = $md2Name
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(3, 9), $md2Url, "", "", $md2Name) | Out-Null
$zhcn.Cells.Item(3, 10).Value = "3791bee1-aff4-47ea-9f14-6c9784f25203.c0eb1211de9929f0bc073d51d8cdf39a57feac62.zh-cn.xlf"
$zhcn.Cells.Item(3, 11).Value = "2016-08-17 12:59:08"

$zhcn.Columns.Item(3).ColumnWidth = 29.95
$zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet, row 2 (0c4fe174...)
# ---------------------------------------------------------------------
$dede.Cells.Item(2, 3).Value = $statusText
$dede.Cells.Item(2, 9).Value = $md1Name
$dede.Hyperlinks.Add($dede.Cells.Item(2, 9), $md1Url, "", "", $md1Name) | Out-Null
$dede.Cells.Item(2, 10).Value = "0c4fe174-d4d5-4219-a86a-488d26943302.c55e58c7ab2c895dc8aced18ebd64001454b9a1e.de-de.xlf"
$dede.Cells.Item(2, 11).Value = "2016-08-17 12:59:16"

# de-de sheet, row 3 (3791bee1...)
$dede.Cells.Item(3, 3).Value = $statusText
$dede.Cells.Item(3, 9).Value = $md2Name
$dede.Hyperlinks.Add($dede.Cells.Item(3, 9), $md2Url, "", "", $md2Name) | Out-Null
$dede.Cells.Item(3, 10).Value = "3791bee1-aff4-47ea-9f14-6c9784f25203.c0eb1211de9929f0bc073d51d8cdf39a57feac62.de-de.xlf"
$dede.Cells.Item(3, 11).Value = "2016-08-17 12:59:16"

$dede.Columns.Item(3).ColumnWidth = 29.95
$dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667

Write-Host "Handback report generated."
